$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mets")
$ws.Columns.Item(4).Delete()
